# The sheet holds a weekly price log for "Bruselas (repollito)" at
# "Vega Modelo de Temuco". A new week's record (2021-10-18 / serial 44487,
# volume 15) is inserted at the top of the data block (row 13), pushing all
# subsequent rows down by one. The last existing row spills into a brand
# new row 62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13; Excel shifts rows 13:61 down to 14:62
# and extends the used range/dimension accordingly.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Vega Modelo de Temuco"
$ws.Range("C13").Value = "La Araucanía"
$ws.Range("D13").Value = 44487
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 100112035
$ws.Range("G13").Value = "Bruselas (repollito)"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 15
$ws.Range("K13").Value = 25000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 25000
$ws.Range("N13").Value = "$/malla 10 kilos"
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 2500
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = "Hortaliza"
